# Update "想去人数" (number of people interested) counts in column F
# for the sheets "展览" and "全部类型", matching the values observed
# in the updated data scrape.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 6649
    6  = 2036
    7  = 1551
    8  = 309
    9  = 1014
    10 = 436
    12 = 5638
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
